# Update the "保險" (Insurance) worksheet (4th sheet) to:
#   1. Fix the header row (it previously duplicated row 2's data instead of
#      real column headers) so it reads:
#      company, name, owner, property_category, category, date,
#      legislator_name, legislator_id, source_file, index
#   2. Append the same metadata columns (property_category, category, date,
#      legislator_name, legislator_id, source_file, index) that already
#      exist on the other sheets (土地/建物/存款) to every data row.
#
# Existing cells elsewhere in the workbook that already hold the exact text
# we need are reused via Copy-to-destination so we don't introduce stray
# number formats (e.g. Excel auto-converting the literal "2013-12-25" text
# into a date) or extra cell styles - the copy just carries over the
# original cell's value/format untouched.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 土地 - source for header labels
$ws3 = $wb.Worksheets.Item(3)   # 存款 - source for repeated metadata values
$ws4 = $wb.Worksheets.Item(4)   # 保險 - sheet being edited

# ---------------------------------------------------------------------
# 1) Header row (row 1)
# ---------------------------------------------------------------------
# B1 becomes the new "company" header (no existing cell holds this text).
$ws4.Range("B1").Value = "company"

# C1 ("name"), D1 ("owner") and the new E1:K1 headers already exist as
# header cells (style s=1, bold) on sheet 1 - reuse them so the shared
# string indexes/styles line up exactly.
$ws1.Range("B1").Copy($ws4.Range("C1")) | Out-Null   # name
$ws1.Range("E1").Copy($ws4.Range("D1")) | Out-Null   # owner
$ws1.Range("I1").Copy($ws4.Range("E1")) | Out-Null   # property_category
$ws1.Range("J1").Copy($ws4.Range("F1")) | Out-Null   # category
$ws1.Range("K1").Copy($ws4.Range("G1")) | Out-Null   # date
$ws1.Range("L1").Copy($ws4.Range("H1")) | Out-Null   # legislator_name
$ws1.Range("M1").Copy($ws4.Range("I1")) | Out-Null   # legislator_id
$ws1.Range("N1").Copy($ws4.Range("J1")) | Out-Null   # source_file
$ws1.Range("O1").Copy($ws4.Range("K1")) | Out-Null   # index

# ---------------------------------------------------------------------
# 2) Data rows (2-5): add property_category, category, date,
#    legislator_name, legislator_id, source_file and index.
# ---------------------------------------------------------------------
$rows = 2,3,4,5
foreach ($r in $rows) {
    $idx = $ws4.Cells.Item($r, 1).Value()

    $ws4.Cells.Item($r, 5).Value = "insurance"

    $ws3.Range("H2").Copy($ws4.Cells.Item($r, 6)) | Out-Null    # category -> normal
    $ws3.Range("I2").Copy($ws4.Cells.Item($r, 7)) | Out-Null    # date -> 2013-12-25
    $ws3.Range("J2").Copy($ws4.Cells.Item($r, 8)) | Out-Null    # legislator_name -> 林明溱
    $ws3.Range("K2").Copy($ws4.Cells.Item($r, 9)) | Out-Null    # legislator_id -> 1706
    $ws3.Range("L2").Copy($ws4.Cells.Item($r, 10)) | Out-Null   # source_file -> tmpfac21

    $ws4.Cells.Item($r, 11).Value = $idx                        # index (mirrors column A)
}
